$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprinklers")

# Sort the sprinkler position table (A10:H58) ascending by column A,
# matching "Validation: changed spreadsheet interval to 1 s to be
# consistent with other sprinkler simulations in the validation suite."
$sortRange = $ws.Range("A10:H58")
$sortKey = $ws.Range("A10:A58")
$sortRange.Sort($sortKey, 1)

# Make the Sprinklers sheet the active tab/sheet and select A10, which is
# where the sort leaves the active cell.
$ws.Activate()
$ws.Range("A10").Select()
